$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Region column (G) updates
$ws.Range("G4:G24").Value = "Norway"
$ws.Range("G25:G35").Value = "Sweden"

# Village column (I) specific updates
$ws.Range("I4").Value = "Sandeid"
$ws.Range("I7").Value = "Belingo"
$ws.Range("I9").Value = "hood"
$ws.Range("I11").Value = "Wiik Village"
$ws.Range("I16").Value = "Cameroon street"
$ws.Range("I25").Value = "Village2"
$ws.Range("I29").Value = "Village5"
$ws.Range("I30").Value = "Village8"
$ws.Range("I32").Value = "Village9"
$ws.Range("I33").Value = "Village10"
$ws.Range("I34").Value = "Village20"

# Touch the bottom-right corner of the used columns to extend the
# worksheet's recorded dimension down to row 1000 (matches the
# observed A1:L1000 dimension in the edited file) without introducing
# any new shared-string/content.
$ws.Cells.Item(1000, 12).NumberFormat = "General"

# Restore the final selection recorded in the edited workbook
$ws.Range("I35").Select()
